# Update "想去人数" (column F) counts on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

# Rows -> new F-column value shared by both "展览" and "全部类型" sheets.
$commonUpdates = @{
    6  = 624
    8  = 340
    9  = 28
    13 = 26
    16 = 110
    17 = 1074
    18 = 1441
    20 = 353
    22 = 91
    24 = 49
    26 = 246
    28 = 295
    29 = 1660
    33 = 613
    35 = 3879
    37 = 451
    39 = 982
    40 = 89
    43 = 88
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $commonUpdates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $commonUpdates[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $commonUpdates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $commonUpdates[$row]
}
# "全部类型" has one extra row updated that "展览" does not.
$ws4.Cells.Item(38, 6).Value = 223
